$wb = $excel.ActiveWorkbook

# --- "BINOMIO DE NEWTON" sheet: fill in the missing rows of the first
# table (rows 6-12, cols B:H) with "OK", and complete the row labels in
# column A (Grado 20 / Grado 50 / Grado 100), mirroring the already
# completed table on the "POLINOMIO" sheet. ---
$wsBin = $wb.Worksheets.Item("BINOMIO DE NEWTON")

$wsBin.Range("A9").Value = "Grado 20"
$wsBin.Range("A10").Value = "Grado 50"
$wsBin.Range("A11").Value = "Grado 100"

# Re-use the formatting already present on B6 (green "OK" fill) for the
# whole B6:H12 block, then stamp every cell with "OK".
$wsBin.Range("B6").Copy()
$okRange = $wsBin.Range("B6:H12")
$okRange.PasteSpecial(-4122)
$okRange.Value = "OK"

# --- Window/selection bookkeeping: update the per-sheet selections and
# finish with "POLINOMIO" as the active sheet/tab (each sheet must be
# activated before its own Range can be selected). ---
$wsPoly = $wb.Worksheets.Item("POLINOMIO")
$wsDesa = $wb.Worksheets.Item("BINOMIO DE NEWTON DESARROLLADO")

$wsBin.Activate()
$wsBin.Range("I11").Select()

$wsDesa.Activate()
$wsDesa.Range("A4").Select()
$wsDesa.Range("J41").Select()

$wsPoly.Activate()
$wsPoly.Range("I23").Select()
